# Insert a new weekly price record for "Camote" (Macroferia Regional de
# Talca - Zapallo) above the current row 449, shifting all the following
# rows down by one (old row 449 -> 450, ... old row 472 -> 473).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 449; Excel shifts rows 449:472 down to
# 450:473 and extends the used range/dimension automatically.
$ws.Rows.Item(449).Insert()

# Populate the newly inserted row 449 with the new record.
$ws.Range("A449").Value = 5
$ws.Range("B449").Value = "Macroferia Regional de Talca"
$ws.Range("C449").Value = "Maule"
$ws.Range("D449").Value = 45147
$ws.Range("E449").Value = 7
$ws.Range("F449").Value = 100112045
$ws.Range("G449").Value = "Zapallo"
$ws.Range("H449").Value = "Camote"
$ws.Range("I449").Value = "1a (guarda)"
$ws.Range("J449").Value = 800
$ws.Range("K449").Value = 350
$ws.Range("L449").Value = 350
$ws.Range("M449").Value = 350
$ws.Range("N449").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O449").Value = "Región del Maule"
$ws.Range("P449").Value = 350
$ws.Range("Q449").Value = 1
$ws.Range("R449").Value = "Hortaliza"
